# Removed Taiwan from List of Countries Due to Lack of Data from the World Bank
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing "Taiwan" in column A and delete the entire row,
# shifting all rows below it up by one.
$found = $ws.Range("A1:A57").Find("Taiwan", [Type]::Missing, [Type]::Missing, 1)
if ($found -ne $null) {
    $ws.Rows($found.Row).Delete()
}

# Restore a plausible selection/view state after the edit.
$ws.Range("E63").Select()
